# Update "想去人数" (number of people interested) figures in the
# 展览 (Exhibition) and 全部类型 (All Types) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 72
$ws1.Range("F3").Value = 803
$ws1.Range("F4").Value = 43
$ws1.Range("F6").Value = 105
$ws1.Range("F8").Value = 4194
$ws1.Range("F9").Value = 92
$ws1.Range("F10").Value = 4890
$ws1.Range("F11").Value = 544
$ws1.Range("F12").Value = 1226

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 72
$ws4.Range("F3").Value = 803
$ws4.Range("F4").Value = 43
$ws4.Range("F6").Value = 105
$ws4.Range("F9").Value = 4194
$ws4.Range("F10").Value = 92
$ws4.Range("F11").Value = 4890
$ws4.Range("F12").Value = 544
$ws4.Range("F13").Value = 1226
